$wb = $excel.ActiveWorkbook

# --- Overview sheet: refresh "Latest HO Xliff Generate Date" for the two
#     files that were just re-handed-back (rows 2 and 4) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-12 16:18:16"
$wsOverview.Range("G4").Value = "2016-08-12 16:18:16"

# --- zh-cn sheet: Priority changed from "ht" to "mt", and the handoff /
#     handback timestamps were refreshed for rows 2 and 4 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-08-12 16:18:09"
$wsZhCn.Range("H4").Value = "2016-08-12 16:18:09"
$wsZhCn.Range("K2").Value = "2016-08-12 16:18:38"
$wsZhCn.Range("K4").Value = "2016-08-12 16:18:38"

# --- de-de sheet: Priority changed from "ht" to "mt" (same shared value as
#     zh-cn), "Correspond Handoff Datetime" shares the same original
#     timestamp as the Overview sheet's date and moves together with it,
#     and the handback timestamp was refreshed for rows 2 and 4 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-08-12 16:18:16"
$wsDeDe.Range("H4").Value = "2016-08-12 16:18:16"
$wsDeDe.Range("K2").Value = "2016-08-12 16:18:47"
$wsDeDe.Range("K4").Value = "2016-08-12 16:18:47"
